$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) / Volume(1h) (E) cells per the Jul 16 2023 cryptos refresh.
$updates = @(
    @{ row=2; d="30.313.62"; e="  -0.10%  " }
    @{ row=3; d="1.931.32"; e="  -0.06%  " }
    @{ row=4; d="0.9991"; e="  -0.25%  " }
    @{ row=5; d="0.7565"; e="  +5.72%  " }
    @{ row=6; d="244.87"; e="  -2.29%  " }
    @{ row=7; d="0.9996"; e="  -0.20%  " }
    @{ row=8; d="0.3188"; e="  -2.50%  " }
    @{ row=9; d="27.61"; e="  -0.13%  " }
    @{ row=10; d="0.06993"; e="  -2.63%  " }
    @{ row=11; d="0.7818"; e="  -2.36%  " }
    @{ row=12; e="  -1.09%  " }
    @{ row=13; d="1.930.46"; e="  -0.11%  " }
    @{ row=14; d="5.358"; e="  -1.17%  " }
    @{ row=15; d="94.38"; e="  -0.19%  " }
    @{ row=16; d="14.43"; e="  -3.00%  " }
    @{ row=17; d="30.315.68"; e="  -0.08%  " }
    @{ row=18; d="252.89"; e="  +0.22%  " }
    @{ row=19; d="0.000007930"; e="  -2.37%  " }
    @{ row=20; d="5.730"; e="  -1.33%  " }
    @{ row=21; d="2.181.18"; e="  -0.15%  " }
    @{ row=22; d="0.9997"; e="  -0.20%  " }
    @{ row=23; d="1.000"; e="  -0.32%  " }
    @{ row=24; d="6.676"; e="  -3.67%  " }
    @{ row=25; d="9.498"; e="  -2.27%  " }
    @{ row=26; d="165.90" }
    @{ row=27; d="18.95"; e="  -1.39%  " }
    @{ row=28; d="0.1331"; e="  +3.23%  " }
    @{ row=29; d="2.215"; e="  -4.99%  " }
    @{ row=30; d="1.362"; e="  -0.26%  " }
    @{ row=31; d="1.512"; e="  -1.91%  " }
    @{ row=32; d="4.381"; e="  -0.93%  " }
    @{ row=33; d="4.120"; e="  -1.89%  " }
    @{ row=34; d="0.05164"; e="  -0.74%  " }
    @{ row=35; d="1.274"; e="  +0.61%  " }
    @{ row=36; d="0.7454"; e="  -0.19%  " }
    @{ row=37; d="2.766"; e="  +0.09%  " }
    @{ row=38; d="0.01949"; e="  -0.54%  " }
    @{ row=39; d="2.792"; e="  -0.24%  " }
    @{ row=40; d="77.89"; e="  -1.42%  " }
    @{ row=41; d="6.409"; e="  -0.80%  " }
    @{ row=42; d="0.4468"; e="  -1.14%  " }
    @{ row=43; d="1.967"; e="  -2.86%  " }
    @{ row=44; d="1.000"; e="  -0.13%  " }
    @{ row=45; d="0.8316"; e="  -0.86%  " }
    @{ row=46; e="  -1.13%  " }
    @{ row=47; d="9.727"; e="  -0.34%  " }
    @{ row=48; e="  +0.63%  " }
    @{ row=49; d="986.39"; e="  +11.79%  " }
    @{ row=50; d="37.15" }
    @{ row=51; d="0.06003"; e="  -1.06%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("d")) {
        $cell = $ws.Range("D" + $u.row)
        # Force text storage so numeric-looking strings (e.g. trailing
        # zeros like "1.000") keep their exact display text instead of
        # being coerced to a trimmed numeric value.
        $cell.NumberFormat = "@"
        $cell.Value = $u.d
        $cell.ClearFormats()
    }
    if ($u.ContainsKey("e")) {
        $ws.Range("E" + $u.row).Value = $u.e
    }
}
